$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2333.3333
$ws.Range("J51").Value = 3000
$ws.Range("L51").Value = 3000
$ws.Range("N51").Value = -3968
$ws.Range("H62").Value = 64175.89
$ws.Range("I62").Value = 94655
$ws.Range("J62").Value = 3217.6667
$ws.Range("K62").Value = 94655
$ws.Range("L62").Value = 3217.6667
$ws.Range("M62").Value = -94031
$ws.Range("N62").Value = -4465.6667
$ws.Range("H65").Value = 64175.89
$ws.Range("I65").Value = 94655
$ws.Range("J65").Value = 3217.6667
$ws.Range("K65").Value = 473275
$ws.Range("L65").Value = 16088.3335
$ws.Range("M65").Value = -470155
$ws.Range("N65").Value = -22328.3335
$ws.Range("H70").Value = 2283.3333
$ws.Range("I70").Value = 1466.6666
$ws.Range("J70").Value = 3100
$ws.Range("K70").Value = 4399.9998
$ws.Range("L70").Value = 9300
$ws.Range("M70").Value = -4129.9998
$ws.Range("N70").Value = -9840
$ws.Range("H73").Value = 2283.3333
$ws.Range("I73").Value = 1466.6666
$ws.Range("J73").Value = 3100
$ws.Range("K73").Value = 4399.9998
$ws.Range("L73").Value = 9300
$ws.Range("M73").Value = -3463.9998
$ws.Range("N73").Value = -11172
$ws.Range("H76").Value = 52383760
$ws.Range("I76").Value = 55002800
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 55002800
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -55002485
$ws.Range("N76").Value = -3630
$ws.Range("H79").Value = 52383760
$ws.Range("I79").Value = 55002800
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 55002800
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -55001708
$ws.Range("N79").Value = -5184
$ws.Range("H98").Value = 1280.75
$ws.Range("I98").Value = 1410.2963
$ws.Range("K98").Value = 1410.2963
$ws.Range("M98").Value = 87.70370000000003
$ws.Range("H111").Value = 736.05884
$ws.Range("I111").Value = 469.5625
$ws.Range("K111").Value = 1408.6875
$ws.Range("M111").Value = 1658.3125
$ws.Range("H113").Value = 2562.611
$ws.Range("I113").Value = 3045.5557
$ws.Range("J113").Value = 2079.6667
$ws.Range("K113").Value = 3045.5557
$ws.Range("L113").Value = 2079.6667
$ws.Range("M113").Value = 208.4443000000001
$ws.Range("N113").Value = -8587.6667
$ws.Range("H116").Value = 3122.5454
$ws.Range("I116").Value = 3256.4285
$ws.Range("J116").Value = 2888.25
$ws.Range("K116").Value = 3256.4285
$ws.Range("L116").Value = 2888.25
$ws.Range("M116").Value = 185.5715
$ws.Range("N116").Value = -9772.25
$ws.Range("H122").Value = 1280.75
$ws.Range("I122").Value = 1410.2963
$ws.Range("K122").Value = 4230.8889
$ws.Range("M122").Value = -1780.8889
$ws.Range("H141").Value = 5015.8335
$ws.Range("I141").Value = 2548.75
$ws.Range("K141").Value = 7646.25
$ws.Range("M141").Value = -2466.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 39999.5
$ws.Range("J92").Value = 39999.5
$ws.Range("L92").Value = 39999.5
$ws.Range("N92").Value = -44991.5
$ws.Range("H132").Value = 2091.9524
$ws.Range("I132").Value = 1275.75
$ws.Range("J132").Value = 2594.2307
$ws.Range("K132").Value = 3827.25
$ws.Range("L132").Value = 7782.6921
$ws.Range("M132").Value = -1297.25
$ws.Range("N132").Value = -12842.6921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 8805.223
$ws.Range("I107").Value = 1575.4615
$ws.Range("J107").Value = 27602.6
$ws.Range("K107").Value = 1575.4615
$ws.Range("L107").Value = 27602.6
$ws.Range("M107").Value = 344.5385000000001
$ws.Range("N107").Value = -31442.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 855.1
$ws.Range("I16").Value = 603.6667
$ws.Range("J16").Value = 962.8570999999999
$ws.Range("K16").Value = 603.6667
$ws.Range("L16").Value = 962.8570999999999
$ws.Range("M16").Value = -316.6667
$ws.Range("N16").Value = -1536.8571
$ws.Range("H113").Value = 855.1
$ws.Range("I113").Value = 603.6667
$ws.Range("J113").Value = 962.8570999999999
$ws.Range("K113").Value = 603.6667
$ws.Range("L113").Value = 962.8570999999999
$ws.Range("M113").Value = 1566.3333
$ws.Range("N113").Value = -5302.8571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 750
$ws.Range("I68").Value = 500
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -689
$ws.Range("N68").Value = -4622
$ws.Range("H71").Value = 750
$ws.Range("I71").Value = 500
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 4500
$ws.Range("L71").Value = 9000
$ws.Range("M71").Value = -444
$ws.Range("N71").Value = -17112
$ws.Range("H131").Value = 778.33
$ws.Range("I131").Value = 289.89474
$ws.Range("J131").Value = 892.90125
$ws.Range("K131").Value = 869.6842200000001
$ws.Range("L131").Value = 2678.70375
$ws.Range("M131").Value = 4170.31578
$ws.Range("N131").Value = -12758.70375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7202
$ws.Range("J70").Value = 9000
$ws.Range("L70").Value = 9000
$ws.Range("N70").Value = -9540
$ws.Range("H73").Value = 7202
$ws.Range("J73").Value = 9000
$ws.Range("L73").Value = 9000
$ws.Range("N73").Value = -10872

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3460
$ws.Range("J68").Value = 2250
$ws.Range("L68").Value = 2250
$ws.Range("N68").Value = -3748
$ws.Range("H71").Value = 3460
$ws.Range("J71").Value = 2250
$ws.Range("L71").Value = 11250
$ws.Range("N71").Value = -18738

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 849.0833
$ws.Range("I107").Value = 867.7368
$ws.Range("J107").Value = 778.2
$ws.Range("K107").Value = 2603.2104
$ws.Range("L107").Value = 2334.6
$ws.Range("M107").Value = -683.2103999999999
$ws.Range("N107").Value = -6174.6
$ws.Range("H113").Value = 252.25
$ws.Range("I113").Value = 257.27274
$ws.Range("K113").Value = 771.81822
$ws.Range("M113").Value = 1398.18178
$ws.Range("H122").Value = 1230.6154
$ws.Range("I122").Value = 1159.8
$ws.Range("K122").Value = 3479.4
$ws.Range("M122").Value = -1029.4
$ws.Range("H140").Value = 31320
$ws.Range("J140").Value = 31320
$ws.Range("L140").Value = 31320
$ws.Range("N140").Value = -41680
